{"js": "// Remove the \"Word version of this document\" list item (and its hyperlink)\n// from the \"Additional resources\" bullet list, per the commit\n// \"Added PDF versions to site\" which drops the old Word-doc download link.\n\nconst searchResults = context.document.body.search(\n  \"Word version of this document\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  const hitParagraphs = searchResults.items[i].paragraphs;\n  hitParagraphs.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < hitParagraphs.items.length; j++) {\n    hitParagraphs.items[j].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Word version of this document\" list item (and its hyperlink)\n# from the \"Additional resources\" bullet list, per the commit\n# \"Added PDF versions to site\" which drops the old Word-doc download link.\n\n$d = $word.ActiveDocument\n$range = $d.Content\n\nif ($range.Find.Execute(\"Word version of this document\")) {\n    $range.Paragraphs(1).Range.Delete()\n}\n"}
